$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4's name changed from "Soohyuk" to "Kim"
$ws.Range("A4").Value = "Kim"

# Fill in the nationality-prediction results (Country/Prob pairs) for each name
$ws.Range("B2").Value = "ID"
$ws.Range("C2").Value = [double]"0.65246074239394103"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = [double]"8.4471059686400101E-3"
$ws.Range("F2").Value = "KW"
$ws.Range("G2").Value = [double]"6.9044722894744601E-3"

$ws.Range("B3").Value = "TH"
$ws.Range("C3").Value = [double]"0.61727555181036498"
$ws.Range("D3").Value = "NO"
$ws.Range("E3").Value = [double]"4.2287434443031499E-2"
$ws.Range("F3").Value = "SG"
$ws.Range("G3").Value = [double]"2.58763910885982E-2"

$ws.Range("B4").Value = "KR"
$ws.Range("C4").Value = [double]"0.52271219397908297"
$ws.Range("D4").Value = "US"
$ws.Range("E4").Value = [double]"2.6445761102455501E-2"
$ws.Range("F4").Value = "VN"
$ws.Range("G4").Value = [double]"2.1798781016200301E-2"

$ws.Range("B5").Value = "HK"
$ws.Range("C5").Value = [double]"0.511905825052546"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = [double]"4.6821068099262397E-2"
$ws.Range("F5").Value = "SG"
$ws.Range("G5").Value = [double]"2.3809548507596898E-2"

$ws.Range("B6").Value = "RU"
$ws.Range("C6").Value = [double]"0.12284491067595101"
$ws.Range("D6").Value = "UA"
$ws.Range("E6").Value = [double]"8.80525791746683E-2"
$ws.Range("F6").Value = "IL"
$ws.Range("G6").Value = [double]"8.5720623601626994E-2"

$ws.Range("B7").Value = "IN"
$ws.Range("C7").Value = [double]"0.29459546763532801"
$ws.Range("D7").Value = "AE"
$ws.Range("E7").Value = [double]"5.5164696570539899E-2"
$ws.Range("F7").Value = "US"
$ws.Range("G7").Value = [double]"3.9912934324786403E-2"

$ws.Range("B8").Value = "US"
$ws.Range("C8").Value = [double]"0.21347775383491099"
$ws.Range("D8").Value = "GB"
$ws.Range("E8").Value = [double]"7.3318437170378503E-2"
$ws.Range("F8").Value = "JM"
$ws.Range("G8").Value = [double]"6.7113294874836693E-2"

# Update the active cell/selection recorded in the sheet view
$ws.Range("J12").Select()
